
# ---------------------------------------------------------------------------
# Kanban con filtro listo
#
# The original workbook has a single sheet "Hoja1" holding a backlog table.
# This edit:
#   1. Renames "Hoja1" -> "Modificacion posterior" (keeps all of its data).
#   2. Inserts a brand-new first sheet "CARGA INICIAL" with a similarly
#      shaped (but filtered/"reset") backlog table, using the same fill
#      colours already present in the workbook's style table.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Sheets: add "CARGA INICIAL" before the existing sheet, rename both ---
$ws1 = $wb.Worksheets.Item(1)
$origName = $ws1.Name

$carga = $wb.Worksheets.Add($ws1)
$carga.Name = "CARGA INICIAL"

$modif = $wb.Worksheets.Item($origName)
$modif.Name = "Modificacion posterior"

# --- 2. Pull formatting donors off "Modificacion posterior" BEFORE we touch
#        its own cosmetics, so "CARGA INICIAL" ends up with the very same
#        fills (theme colours), just reused from the existing style table. ---

# Row 1 styles (merged-header look): fillId3+center, hyperlink-blank,
# fillId4+center, fillId5+center
$modif.Range("B1").Copy()
$carga.Range("B1:C1").PasteSpecial(-4122)

$modif.Range("D3").Copy()
$carga.Range("D1").PasteSpecial(-4122)

$modif.Range("G1").Copy()
$carga.Range("E1:G1").PasteSpecial(-4122)

$modif.Range("K1").Copy()
$carga.Range("H1:M1").PasteSpecial(-4122)

# Row 2 styles (column headers): fillId3, fillId2, fillId4, fillId5 (no align)
$modif.Range("B2").Copy()
$carga.Range("B2:C2").PasteSpecial(-4122)

$modif.Range("F2").Copy()
$carga.Range("D2").PasteSpecial(-4122)

$modif.Range("G2").Copy()
$carga.Range("E2:G2").PasteSpecial(-4122)

$modif.Range("K2").Copy()
$carga.Range("H2:M2").PasteSpecial(-4122)

# Data rows: column D (assignee) keeps the fillId5 "no align" look throughout
$modif.Range("K2").Copy()
$carga.Range("D3:D12").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 3. Merge the header cells ---
$carga.Range("B1:C1").Merge()
$carga.Range("H1:M1").Merge()

# --- 4. Values ---
$carga.Range("B1").Value = "Necesarias"
$carga.Range("E1").Value = "Para despues"
$carga.Range("H1").Value = "Atributos Adicionales (etiquetas en gitlab)"

$carga.Range("A2").Value = "Id (Excel)"
$carga.Range("B2").Value = "Nombre 1"
$carga.Range("C2").Value = "Descripción 2"
$carga.Range("D2").Value = "Asignación (login) 4"
$carga.Range("E2").Value = "iterasion 4"
$carga.Range("F2").Value = "prioridad 5"
$carga.Range("G2").Value = "Estatus 6"
$carga.Range("H2").Value = "Puntos 7"
$carga.Range("I2").Value = "Talla"

$lorem = "Lorem ipsum dolor sit amet consectetur adipisicing elit. Accusamus vel architecto magnam culpa asperiores voluptates omnis maxime neque quas aliquam qui obcaecati itaque exercitationem, enim cum consectetu"

$rows = @(
  @(1,  "Historia de Usuario 1",         "admin, dev1", "Sin Iteracion", "Alta",  "Por Asignar", 10, "ch"),
  @(2,  "Historia de Usuario 2",         "admin, dev1", "Sin Iteracion", "Alta",  "Por Asignar", 20, "G"),
  @(3,  "Historia de Usuario 3",         "admin",       "Sin Iteracion", "Baja",  "Por Asignar", 30, "M"),
  @(4,  "Historia de Usuario 4",         "admin",       "Sin Iteracion", "Media", "Por Asignar", 10, "ch"),
  @(5,  "Historia de Usuario 5",         "admin",       "Sin Iteracion", "Alta",  "Por Asignar", 20, "M"),
  @(6,  "Historia de Usuario 6",         "admin",       "Sin Iteracion", "Alta",  "Por Asignar", 30, "ch"),
  @(7,  "Historia de Usuario 7",         "admin",       "Sin Iteracion", "Baja",  "Por Asignar", 70, "M"),
  @(8,  "Historia de Usuario 8",         "admin",       "Sin Iteracion", "Media", "Por Asignar", 10, "ch"),
  @(9,  "Historia de Usuario 9",         "admin",       "Sin Iteracion", "Alta",  "Por Asignar", 20, "G"),
  @(10, "Test Historia de Usuario 10",   "admin",       "Por Asignar",   "Alta",  "Por Asignar", 30, "ch")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 3
    $row = $rows[$i]
    $carga.Cells.Item($r, 1).Value = $row[0]
    $carga.Cells.Item($r, 2).Value = $row[1]
    $carga.Cells.Item($r, 3).Value = $lorem
    $carga.Cells.Item($r, 4).Value = $row[2]
    $carga.Cells.Item($r, 5).Value = $row[3]
    $carga.Cells.Item($r, 6).Value = $row[4]
    $carga.Cells.Item($r, 7).Value = $row[5]
    $carga.Cells.Item($r, 8).Value = $row[6]
    $carga.Cells.Item($r, 9).Value = $row[7]
}

# --- 5. Column widths (characters) ---
$carga.Columns.Item(2).ColumnWidth = 28.6
$carga.Columns.Item(4).ColumnWidth = 19.6
$carga.Range("E1:G1").ColumnWidth = 10.8

# --- 6. View state on both sheets ---
$modif.Activate()
$modif.Range("B58").Select()

$carga.Activate()
$carga.Range("C21").Select()
